# Apply the three changes described by the commit diff:
#   1. Table: switch to fixed table layout (<w:tblLayout w:type="fixed"/>)
#   2. Styles: introduce a new "Abstract Title" paragraph style
#   3. Styles: tighten the "Abstract" style's space-before (300 -> 100)

$d = $word.ActiveDocument

# --- 1. Table layout -------------------------------------------------
$tbl = $d.Tables(1)
$tbl.AllowAutoFit = $false

# --- 2. New "Abstract Title" style ------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 3. Tighten "Abstract" style spacing ------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

Write-Output "Applied tblLayout fixed, added Abstract Title style, tightened Abstract spacing."
